$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.242630004882812
$ws.Range("B1").Value = 3.850899934768677
$ws.Range("C1").Value = 3.570510149002075
$ws.Range("D1").Value = 3.587673187255859
$ws.Range("E1").Value = 1.092154383659363
